$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Expand the WBS table by two columns (it currently ends at column G / Remark).
# The new columns ("Labour Count", "Productivity Ref") are inserted just before
# the trailing "Remark" column, so first grow the table to A1:I48 (Excel appends
# the two new columns at the end), then shift the header text one-by-one so that
# "Remark" ends up last again.
$lo.Resize($ws.Range("A1:I48"))

$ws.Range("G1").Value = "Labour Count"
$ws.Range("H1").Value = "Productivity Ref"
$ws.Range("I1").Value = "Remark"

# Copy the header cell formatting (bold/centered style) from the existing
# "Equation $v" header onto the three headers that moved/were created so they
# keep looking like proper table headers.
$ws.Range("F1").Copy()
$ws.Range("G1:I1").PasteSpecial(-4122)

# Match column widths: the two new columns inherit the width of the column
# immediately to their left ("Equation $v"), while "Remark" keeps the width it
# had before the insert, and the trailing formatted-but-empty column shifts
# two slots to the right along with its width.
$ws.Range("G:H").ColumnWidth = 21.33
$ws.Range("I:I").ColumnWidth = 12.83
$ws.Range("J:J").ColumnWidth = 19.83

# The workbook was re-saved from a different folder on disk.
$wb.Path = "/Users/hazem/Downloads/"

[void]$ws.Range("F10").Select()
